$wb = $excel.ActiveWorkbook

# --- "BLS Data Series" tab holds the year-by-month unemployment data table ---
$data = $wb.Worksheets.Item("BLS Data Series")

# Drop the 2000-2002 rows (no Asian unemployment data before 2003, per commit
# message) so the table starts at 2003 and shifts everyone up three rows.
$data.Range("A2:A4").EntireRow.Delete()

# Restore the view/selection state recorded in the saved file: frozen header
# pane, with the live selection anchored at A2 covering A2:XFD4.
$data.Activate()
$data.Range("A2:XFD4").Select()

# --- "Sheet1" tab holds the BLS series metadata/description block ---
$meta = $wb.Worksheets.Item("Sheet1")
$meta.Rows.Item(7).RowHeight = 28
$meta.Rows.Item(8).RowHeight = 28
